$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the columns that no longer exist in the new schema (H:L) — entire
# column delete so the used range / dimension shrinks to A:G.
$ws.Range("H1:L5").EntireColumn.Delete()

# Three more data rows are needed (rows 6,7,8 beyond the original A1:G5
# extent). Prime their column-A formatting (bordered/bold/centred style
# used by every "glycan" cell) by copying it down from the existing A5
# before any values are written, so the new cells pick up style id 1
# instead of falling back to the default.
$ws.Range("A5").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "glycan"
$ws.Range("B1").Value = "binding_score"
$ws.Range("C1").Value = "monosaccharides"
$ws.Range("D1").Value = "motifs"
$ws.Range("E1").Value = "sasa"
$ws.Range("F1").Value = "flexibility"
$ws.Range("G1").Value = "has_multi_node_motifs"

# --- Data rows ----------------------------------------------------------
# glycan, binding_score, monosaccharides, motifs, sasa, flexibility, has_multi_node_motifs
$rows = @(
    @("Fuc(a1-3)[Gal(b1-4)]GlcNAc(b1-3)[Fuc(a1-3)[Gal(b1-4)]GlcNAc(b1-6)]Gal(b1-4)Glc", -0.8040186117466136, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 3.980328598306347, 4.477033149590943, $true),
    @("Fuc(a1-3)[Gal(b1-4)]GlcNAc(b1-3)[GlcNAc(b1-6)]Gal(b1-4)Glc", -0.8040186117466136, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 4.145843025381176, 2.375316762015688, $true),
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc", -0.016104048483271, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.338285572580087, 0.91117855161729, $true),
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc", -0.2726764879960648, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.27278254643194, 2.180924532303609, $true),
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc", -0.1662041962938344, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.30612914780972, 3.097700174032828, $true),
    @("Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)GlcNAc", 0.4407208778426784, "['Gal(b1-4)', 'GlcNAc(b1-3)']", "['Gal(b1-4)GlcNAc(b1-3)']", 5.252147263686476, 1.839554809126105, $true),
    @("GlcNAc(b1-4)GlcNAc(b1-4)GlcNAc", -0.1011239779733945, "['GlcNAc(b1-4)', 'GlcNAc(b1-4)']", "['GlcNAc(b1-4)GlcNAc']", 5.702079978569953, 0.6713220512263312, $true)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
